# Applies the "Add data for 2021-11-27" update:
#  - Renames the sheet / updates the running title from "...November 18" to "...November 19"
#  - Updates the carjacking counts for several neighborhood/month cells to reflect the
#    newly added day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name and running-total title -------------------------------------------------
$wb.Sheets.Item(1).Name = "Through 2021-11-19"
$ws.Range("B1").Value = "November 2021 (through November 19)"

# --- Updated / newly populated counts ---------------------------------------------------
$ws.Range("M2").Value   = 11
$ws.Range("AT2").Value  = 3

$ws.Range("M4").Value   = 6
$ws.Range("BE4").Value  = 2

$ws.Range("B10").Value  = 1
$ws.Range("M10").Value  = 4

$ws.Range("M13").Value  = 2

$ws.Range("M15").Value  = 3
$ws.Range("BP15").Value = 1

$ws.Range("B16").Value  = 3

$ws.Range("AT17").Value = 5

$ws.Range("B23").Value  = 3

$ws.Range("AT24").Value = 1

$ws.Range("X26").Value  = 1

$ws.Range("B31").Value  = 5
$ws.Range("M31").Value  = 4

$ws.Range("M38").Value  = 1

$ws.Range("B39").Value  = 3

$ws.Range("M49").Value  = 1
$ws.Range("X49").Value  = 1

$ws.Range("B56").Value  = 1

$ws.Range("BE68").Value = 1

$ws.Range("B84").Value  = 3
$ws.Range("AT84").Value = 1
